# Apply scheduled-runner value updates to Leviathan_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR tables, per diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value2 = 481.6
$ws.Range("I33").Value2 = 477
$ws.Range("K33").Value2 = 477
$ws.Range("M33").Value2 = -248
# Row 70
$ws.Range("H70").Value2 = 3427.0312
$ws.Range("I70").Value2 = 3280.2778
$ws.Range("K70").Value2 = 9840.8334
$ws.Range("M70").Value2 = -9570.8334
# Row 73
$ws.Range("H73").Value2 = 3427.0312
$ws.Range("I73").Value2 = 3280.2778
$ws.Range("K73").Value2 = 9840.8334
$ws.Range("M73").Value2 = -8904.8334
# Row 98
$ws.Range("H98").Value2 = 1460.9412
$ws.Range("I98").Value2 = 774.2143
$ws.Range("K98").Value2 = 774.2143
$ws.Range("M98").Value2 = 723.7857
# Row 100
$ws.Range("H100").Value2 = 3919.0435
$ws.Range("I100").Value2 = 1796.2778
$ws.Range("J100").Value2 = 11561
$ws.Range("K100").Value2 = 1796.2778
$ws.Range("L100").Value2 = 11561
$ws.Range("M100").Value2 = -1255.2778
$ws.Range("N100").Value2 = -12643
# Row 112
$ws.Range("H112").Value2 = 1396.238
$ws.Range("J112").Value2 = 1454.4375
$ws.Range("L112").Value2 = 4363.3125
$ws.Range("N112").Value2 = -6579.3125
# Row 122
$ws.Range("H122").Value2 = 1460.9412
$ws.Range("I122").Value2 = 774.2143
$ws.Range("K122").Value2 = 2322.6429
$ws.Range("M122").Value2 = 127.3571000000002
# Row 132
$ws.Range("H132").Value2 = 3843.375
$ws.Range("I132").Value2 = 1213
$ws.Range("K132").Value2 = 3639
$ws.Range("M132").Value2 = -1109
# Row 138
$ws.Range("H138").Value2 = 2711.366
$ws.Range("J138").Value2 = 3164.2964
$ws.Range("L138").Value2 = 9492.889200000001
$ws.Range("N138").Value2 = -19772.8892

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value2 = 2474.7
$ws.Range("I2").Value2 = 2265.074
$ws.Range("K2").Value2 = 2265.074
$ws.Range("M2").Value2 = -2152.074
# Row 61
$ws.Range("H61").Value2 = 1627.2759
$ws.Range("I61").Value2 = 1541.6305
$ws.Range("K61").Value2 = 1541.6305
$ws.Range("M61").Value2 = -1329.6305
# Row 110
$ws.Range("H110").Value2 = 2559.5
$ws.Range("I110").Value2 = 1261.9375
$ws.Range("K110").Value2 = 1261.9375
$ws.Range("M110").Value2 = 783.0625
# Row 116
$ws.Range("H116").Value2 = 2474.7
$ws.Range("I116").Value2 = 2265.074
$ws.Range("K116").Value2 = 2265.074
$ws.Range("M116").Value2 = 28.92599999999993
# Row 136
$ws.Range("H136").Value2 = 1627.2759
$ws.Range("I136").Value2 = 1541.6305
$ws.Range("K136").Value2 = 4624.8915
$ws.Range("M136").Value2 = -2074.8915

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value2 = 2474.7
$ws.Range("I3").Value2 = 2265.074
$ws.Range("K3").Value2 = 2265.074
$ws.Range("M3").Value2 = -2151.074

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 26199.62
$ws.Range("I31").Value2 = 32379.455
$ws.Range("K31").Value2 = 32379.455
$ws.Range("M31").Value2 = -32084.455
# Row 34
$ws.Range("H34").Value2 = 26199.62
$ws.Range("I34").Value2 = 32379.455
$ws.Range("K34").Value2 = 32379.455
$ws.Range("M34").Value2 = -32177.455
# Row 86
$ws.Range("H86").Value2 = 3296.8333
$ws.Range("J86").Value2 = 3336.2856
$ws.Range("L86").Value2 = 3336.2856
$ws.Range("N86").Value2 = -5582.2856
# Row 89
$ws.Range("H89").Value2 = 3296.8333
$ws.Range("J89").Value2 = 3336.2856
$ws.Range("L89").Value2 = 16681.428
$ws.Range("N89").Value2 = -27913.428
# Row 122
$ws.Range("H122").Value2 = 55811.95
$ws.Range("I122").Value2 = 74362.86
$ws.Range("K122").Value2 = 223088.58
$ws.Range("M122").Value2 = -220638.58
# Row 134
$ws.Range("H134").Value2 = 2411.6086
$ws.Range("I134").Value2 = 1770.6666
$ws.Range("K134").Value2 = 5311.9998
$ws.Range("M134").Value2 = -2776.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value2 = 528
$ws.Range("I10").Value2 = 455.33334
$ws.Range("J10").Value2 = 600.6667
$ws.Range("K10").Value2 = 1366.00002
$ws.Range("L10").Value2 = 1802.0001
$ws.Range("M10").Value2 = -1227.00002
$ws.Range("N10").Value2 = -2080.0001
# Row 19
$ws.Range("H19").Value2 = 549.25
$ws.Range("J19").Value2 = 549.25
$ws.Range("L19").Value2 = 1647.75
$ws.Range("N19").Value2 = -1995.75
# Row 26
$ws.Range("H26").Value2 = 191.64285
$ws.Range("I26").Value2 = 63.6
$ws.Range("J26").Value2 = 262.77777
$ws.Range("K26").Value2 = 190.8
$ws.Range("L26").Value2 = 788.33331
$ws.Range("M26").Value2 = 97.19999999999999
$ws.Range("N26").Value2 = -1364.33331
# Row 37
$ws.Range("H37").Value2 = 142886990
$ws.Range("J37").Value2 = 142886990
$ws.Range("L37").Value2 = 428660970
$ws.Range("N37").Value2 = -428661194
# Row 41
$ws.Range("H41").Value2 = 356.55
$ws.Range("I41").Value2 = 213.1
$ws.Range("J41").Value2 = 500
$ws.Range("K41").Value2 = 639.3
$ws.Range("L41").Value2 = 1500
$ws.Range("M41").Value2 = -301.3
$ws.Range("N41").Value2 = -2176
# Row 68
$ws.Range("H68").Value2 = 1597.8
$ws.Range("J68").Value2 = 1749.5
$ws.Range("L68").Value2 = 5248.5
$ws.Range("N68").Value2 = -6870.5
# Row 71
$ws.Range("H71").Value2 = 1597.8
$ws.Range("J71").Value2 = 1749.5
$ws.Range("L71").Value2 = 15745.5
$ws.Range("N71").Value2 = -23857.5
# Row 81
$ws.Range("H81").Value2 = 100009200
$ws.Range("I81").Value2 = 7997.5
$ws.Range("J81").Value2 = 166676670
$ws.Range("K81").Value2 = 23992.5
$ws.Range("L81").Value2 = 500030010
$ws.Range("M81").Value2 = -22869.5
$ws.Range("N81").Value2 = -500032256
# Row 84
$ws.Range("H84").Value2 = 100009200
$ws.Range("I84").Value2 = 7997.5
$ws.Range("J84").Value2 = 166676670
$ws.Range("K84").Value2 = 71977.5
$ws.Range("L84").Value2 = 1500090030
$ws.Range("M84").Value2 = -66361.5
$ws.Range("N84").Value2 = -1500101262
# Row 122
$ws.Range("H122").Value2 = 617.3333
$ws.Range("J122").Value2 = 579.6
$ws.Range("L122").Value2 = 5216.400000000001
$ws.Range("N122").Value2 = -10116.4
# Row 131
$ws.Range("H131").Value2 = 49345.477
$ws.Range("J131").Value2 = 2298.125
$ws.Range("L131").Value2 = 6894.375
$ws.Range("N131").Value2 = -16974.375
# Row 137
$ws.Range("H137").Value2 = 3901.6667
$ws.Range("J137").Value2 = 4470
$ws.Range("L137").Value2 = 13410
$ws.Range("N137").Value2 = -23610

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value2 = 18872.807
$ws.Range("I7").Value2 = 31139.889
$ws.Range("J7").Value2 = 13854.454
$ws.Range("K7").Value2 = 31139.889
$ws.Range("L7").Value2 = 13854.454
$ws.Range("M7").Value2 = -31027.889
$ws.Range("N7").Value2 = -14078.454
# Row 22
$ws.Range("H22").Value2 = 797.86957
$ws.Range("I22").Value2 = 649.4286
$ws.Range("J22").Value2 = 862.8125
$ws.Range("K22").Value2 = 649.4286
$ws.Range("L22").Value2 = 862.8125
$ws.Range("M22").Value2 = -354.4286
$ws.Range("N22").Value2 = -1452.8125
# Row 27
$ws.Range("H27").Value2 = 797.86957
$ws.Range("I27").Value2 = 649.4286
$ws.Range("J27").Value2 = 862.8125
$ws.Range("K27").Value2 = 649.4286
$ws.Range("L27").Value2 = 862.8125
$ws.Range("M27").Value2 = -542.4286
$ws.Range("N27").Value2 = -1076.8125
# Row 126
$ws.Range("H126").Value2 = 18872.807
$ws.Range("I126").Value2 = 31139.889
$ws.Range("J126").Value2 = 13854.454
$ws.Range("K126").Value2 = 93419.667
$ws.Range("L126").Value2 = 41563.362
$ws.Range("M126").Value2 = -90949.667
$ws.Range("N126").Value2 = -46503.362

$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value2 = 19018.572
$ws.Range("J70").Value2 = 19018.572
$ws.Range("L70").Value2 = 19018.572
$ws.Range("N70").Value2 = -19648.572
# Row 73
$ws.Range("H73").Value2 = 19018.572
$ws.Range("J73").Value2 = 19018.572
$ws.Range("L73").Value2 = 19018.572
$ws.Range("N73").Value2 = -21202.572
# Row 94
$ws.Range("H94").Value2 = 25219.666
$ws.Range("I94").Value2 = 0
$ws.Range("J94").Value2 = 25219.666
$ws.Range("K94").Value2 = 0
$ws.Range("L94").Value2 = 25219.666
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value2 = -27021.666
# Row 122
$ws.Range("H122").Value2 = 1625.8718
$ws.Range("I122").Value2 = 1562.2333
$ws.Range("K122").Value2 = 4686.699900000001
$ws.Range("M122").Value2 = -2236.699900000001
# Row 126
$ws.Range("H126").Value2 = 1880.2222
$ws.Range("I126").Value2 = 1880.2222
$ws.Range("K126").Value2 = 5640.6666
$ws.Range("M126").Value2 = -3170.6666

